$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells G1 and H1, copying the bold header style from F1
$ws.Range("F1").Copy($ws.Range("G1:H1"))
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"

# Update recalculated metrics (B:D) for each row
$ws.Range("B2").Value = 0.5127686906183955
$ws.Range("C2").Value = 0.9897896864432949
$ws.Range("D2").Value = 0.5810528266363498

$ws.Range("B3").Value = 0.2424375111208306
$ws.Range("C3").Value = 0.9952615527482974
$ws.Range("D3").Value = 0.3822725189317837

$ws.Range("B4").Value = 0.3424764921622552
$ws.Range("C4").Value = 0.9934115138357396
$ws.Range("D4").Value = 0.4643252469896818

$ws.Range("B5").Value = 0.4120572361557082
$ws.Range("C5").Value = 0.9918747954596656
$ws.Range("D5").Value = 0.5098666015399509

# Populate new Elapsed Time / CPU columns (G:H) for each row
$ws.Range("G2").Value = 0.4788041146331428
$ws.Range("H2").Value = 0.997

$ws.Range("G3").Value = 0.4788041146331428
$ws.Range("H3").Value = 0.997

$ws.Range("G4").Value = 0.4788041146331428
$ws.Range("H4").Value = 0.997

$ws.Range("G5").Value = 0.4788041146331428
$ws.Range("H5").Value = 0.997
